# Revise the "Group" roster: Rehab and Songtao's group now also includes
# Jiahui, so update the membership text accordingly, then leave the
# "Group" sheet active with B3 selected (mirrors the author's final
# on-screen state after editing the roster).

$wb = $excel.ActiveWorkbook

$wsGroup = $wb.Worksheets.Item("Group")
$wsGroup.Range("B2").Value = "Rehab, Songtao, Jiahui"

$wsGroup.Activate()
[void]$wsGroup.Range("B3").Select()
